# Changes in the Doctor Tests
# Populate the "Actual Result" (column N) values for the CreateDoctor sheet,
# rows 2-8, mirroring / introducing the recorded actual test outcomes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateDoctor")

$ws.Cells.Item(2, 14).Value = "Email ID Already used, please try to login with credentials or use another Email ID"
$ws.Cells.Item(3, 14).Value = "Phone Number Already used, please try to login with credentials or use another Phone Number"
$ws.Cells.Item(4, 14).Value = "Degree is Required"
$ws.Cells.Item(5, 14).Value = "Invalid email address"
$ws.Cells.Item(6, 14).Value = "Please enter a Valid Mobile Number"
$ws.Cells.Item(7, 14).Value = "Name is Required"
$ws.Cells.Item(8, 14).Value = "Experience is Required"
